# Grey out the "Game Design:" sub-section under "Others" — same treatment
# already used elsewhere in the doc for finished/obsolete sections:
# font color = theme color "Background 1" shaded to A6 (i.e. w:color
# w:val="A6A6A6" w:themeColor="background1" w:themeShade="A6").
#
# Packed WdColor value for that exact theme color + shade combo, read
# straight off an existing run in this document that already carries
# <w:color w:val="A6A6A6" w:themeColor="background1" w:themeShade="A6"/>.
$greyThemeColor = -603937025

$d = $word.ActiveDocument

# Exact paragraph texts (as they read via Range.Text, i.e. without the
# trailing paragraph mark) that must turn grey. The list-item about
# instantiating the purchased plot ("Sau khi ng dùng ...") is left as-is —
# it was not part of the requested change.
$targetTexts = @(
    "Game Design:",
    "Bản đồ bao gồm các ô đất và các con đường đc sắp xếp theo hình bàn cờ",
    "Các con đường và ô đất chưa mua sẽ chỉ là background (ô đất thì có collider)",
    "Các con đường xung quanh các ô đất -> khi bấm vào 1 điểm trên con đường -> xác định ô đường của điểm nhấn -> di chuyển nhân vật chính đến ô đó = Shortest Path Algorithm"
)

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs.Item($i)
    # Range.Text includes the trailing paragraph-mark character (CR, and
    # sometimes a BEL for cell marks) — strip it before comparing.
    $text = $para.Range.Text.TrimEnd([char]13, [char]7)

    foreach ($target in $targetTexts) {
        if ($text -eq $target) {
            # Colour the whole paragraph range (run text + the paragraph
            # mark itself), matching Word's own behaviour when you select
            # the paragraph and apply a font colour from the UI.
            $para.Range.Font.Color = $greyThemeColor
        }
    }
}
